# Insert a new data row at row 399 (pushing the existing rows 399-423 down
# to 400-424, which is how the sheet grows from A1:T423 to A1:T424), then
# populate the newly-inserted row with the new "Naranja / Valencia" record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(399).Insert()

$ws.Range("A399").Value = 5
$ws.Range("B399").Value = "Macroferia Regional de Talca"
$ws.Range("C399").Value = "Maule"
$ws.Range("D399").Value = 44585
$ws.Range("E399").Value = 7
$ws.Range("F399").Value = "Fruta"
$ws.Range("G399").Value = 100102
$ws.Range("H399").Value = "Cítricos"
$ws.Range("I399").Value = 100102005
$ws.Range("J399").Value = "Naranja"
$ws.Range("K399").Value = "Valencia"
$ws.Range("L399").Value = "Primera"
$ws.Range("M399").Value = 300
$ws.Range("N399").Value = 10000
$ws.Range("O399").Value = 10000
$ws.Range("P399").Value = 10000
$ws.Range("Q399").Value = "$/bandeja 15 kilos granel"
$ws.Range("R399").Value = "Cabildo"
$ws.Range("S399").Value = 667
$ws.Range("T399").Value = 15
